$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "CB23002"
$ws.Range("B7").Value = 12345678

$ws.Range("B7").Select()
